$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Wed Sep  4 11:21:07 UTC 2024 with GitHub Actions

$ws.Range("D2").Value = '56.366.71'
$ws.Range("E2").Value = '  -4.31%  '

$ws.Range("D3").Value = '2.391.50'
$ws.Range("E3").Value = '  -4.46%  '

$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '500.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.35%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.89%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.551'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.16%  '

$ws.Range("D9").Value = '2.388.49'
$ws.Range("E9").Value = '  -4.71%  '

$ws.Range("E10").Value = '  -4.09%  '

$ws.Range("E11").Value = '  -1.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.318'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.76%  '

$ws.Range("E13").Value = '  -11.30%  '

$ws.Range("D14").Value = '2.815.94'
$ws.Range("E14").Value = '  -4.42%  '

$ws.Range("D15").Value = '57.070.46'
$ws.Range("E15").Value = '  -2.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.44'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.06%  '

$ws.Range("E17").Value = '  -3.84%  '

$ws.Range("D18").Value = '2.407.23'
$ws.Range("E18").Value = '  -3.77%  '

$ws.Range("E19").Value = '  -5.04%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '310.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.28%  '

$ws.Range("E21").Value = '  -5.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.53%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.54%  '

$ws.Range("D26").Value = '2.495.80'
$ws.Range("E26").Value = '  -4.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.371'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -9.30%  '

$ws.Range("E28").Value = '  -6.20%  '

$ws.Range("E29").Value = '  -3.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '174.01'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.10%  '

$ws.Range("E31").Value = '  -4.71%  '

$ws.Range("D32").Value = '0.0₃0704'
$ws.Range("E32").Value = '  -6.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.04%  '

$ws.Range("E35").Value = '  -7.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.997'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.69'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.21'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.04%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.73'
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.62%  '

$ws.Range("E41").Value = '  -6.45%  '

$ws.Range("E42").Value = '  -7.03%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '128.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.68%  '

$ws.Range("E44").Value = '  -4.10%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.76'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.42%  '

$ws.Range("E46").Value = '  -3.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '252.61'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0895'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.23%  '

$ws.Range("E49").Value = '  -5.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.66%  '

$ws.Range("E51").Value = '  -5.13%  '
